$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet (also updates sheet name in workbook.xml)
$ws.Name = "Through 2022-09-17"

# Update the header label for the "2022" column (shared string / cell I1)
$ws.Range("I1").Value = "2022 (through 09-17)"

# Update data values
$ws.Range("I10").Value = 82
$ws.Range("I14").Value = 1217
